# PE_C32017.xlsx - "Finalizado con todo de Programacion Estructurada"
#
# The author added an "X" column score of 30 (a fixed bonus/penalty) to most
# rows of the gradebook, which automatically ripples into the Y (total)
# column because Y is a formula (P+...+U*V+W+X+Z). Two rows (24 and 89) had
# their previous X value of 10 replaced by -1 instead. Row 61 additionally
# got its U value bumped from 16 to 19. Finally the sheet's frozen-pane
# window was scrolled down/right a bit and the active selection moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column X (and the one U tweak) ---------------------------------------
# Rows that gained a new X value of 30
$rowsWith30 = 20,21,30,31,32,33,38,42,47,50,53,57,59,61,62,68,70,71,84,90,91
foreach ($r in $rowsWith30) {
    $ws.Range("X$r").Value = 30
}

# Rows whose X value changed from 10 to -1
$rowsWithMinus1 = 24,89
foreach ($r in $rowsWithMinus1) {
    $ws.Range("X$r").Value = -1
}

# Row 61 also had U61 changed from 16 to 19
$ws.Range("U61").Value = 19

# The Y column holds shared formulas that total each row, so setting the
# cells above is enough for Excel to recompute Y20, Y21, Y24, ... on its own.

# --- Sheet view: scroll / selection ----------------------------------------
# Keep the existing freeze (10 header rows, no column split), just move the
# window and the active selection the way the author left them.
$aw = $excel.ActiveWindow
$aw.FreezePanes = $false
$ws.Range("A11").Select()
$aw.FreezePanes = $true
$ws.Range("X38").Select()
